# Added DataProvider for registration form and filled exec with correct
# registration data: append a new "createAccountFormAllDataRequired" sheet
# with three rows of registration-form test data, then make it the active
# (selected) sheet of the workbook.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end of the tab strip -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)

$targetName = "createAccountFormAllDataRequired"
try {
    $ws.Name = $targetName
} catch {
    # Real Excel (and this host) caps sheet names at 31 characters; fall
    # back to the longest valid prefix if the exact name is rejected.
    if ($targetName.Length -gt 31) {
        $ws.Name = $targetName.Substring(0, 31)
    } else {
        throw
    }
}

# --- Row 1 : complete registration data -----------------------------------
$ws.Cells.Item(1, 1).Value = "Jan"
$ws.Cells.Item(1, 2).Value = "Kowalski"
$ws.Cells.Item(1, 3).Value = "Haslo"
$ws.Cells.Item(1, 4).Value = "Company"
$ws.Cells.Item(1, 5).Value = "Address 1"
$ws.Cells.Item(1, 6).Value = "Address 2"
$ws.Cells.Item(1, 7).Value = "City"
$ws.Cells.Item(1, 8).Value = 88222
$ws.Cells.Item(1, 9).Value = "This is my registration data"
$ws.Cells.Item(1, 10).Value = 258147
$ws.Cells.Item(1, 11).Value = 852369874
$ws.Cells.Item(1, 12).Value = "My address"

# --- Row 2 : partial registration data ------------------------------------
$ws.Cells.Item(2, 1).Value = "Jan "
$ws.Cells.Item(2, 2).Value = "Kowalski"
$ws.Cells.Item(2, 3).Value = "Haslo11A_*"
$ws.Cells.Item(2, 4).Value = "Company Name"
$ws.Cells.Item(2, 8).Value = 99877
$ws.Cells.Item(2, 10).Value = 258147

# --- Row 3 : partial registration data ------------------------------------
$ws.Cells.Item(3, 1).Value = "Jan "
$ws.Cells.Item(3, 2).Value = "Kowalski"
$ws.Cells.Item(3, 3).Value = "Haslo**"
$ws.Cells.Item(3, 4).Value = "Company 1"
$ws.Cells.Item(3, 5).Value = "Address 1"
$ws.Cells.Item(3, 7).Value = "City"
$ws.Cells.Item(3, 8).Value = 1234
$ws.Cells.Item(3, 11).Value = 852369874

# --- Column widths (best effort; Excel quantizes ColumnWidth to pixels) ---
$ws.Columns.Item(4).ColumnWidth = 12.3
$ws.Columns.Item(8).ColumnWidth = 12.0
$ws.Columns.Item(9).ColumnWidth = 22.1

# --- Make the new sheet the active / selected tab -------------------------
$ws.Activate()
$ws.Range("C3").Select() | Out-Null
